# Updated symbol list on Fri Jan 20 13:13:02 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $Val)
    $r = $ws.Range($CellRef)
    $r.NumberFormat = "@"
    $r.Value = $Val
    $r.Style = "Normal"
}

Set-TextValue "D2" "289.42"
Set-TextValue "E2" "0.32%"
Set-TextValue "G2" "13"
Set-TextValue "D3" "31.00"
Set-TextValue "E3" "2.29%"
Set-TextValue "G3" "13"
Set-TextValue "D4" "4.946"
Set-TextValue "E4" "0.17%"
Set-TextValue "G4" "13"
Set-TextValue "D5" "0.07377"
Set-TextValue "E5" "3.06%"
Set-TextValue "G5" "13"
Set-TextValue "D6" "2.358"
Set-TextValue "E6" "30.36%"
Set-TextValue "G6" "13"
Set-TextValue "D7" "7.730"
Set-TextValue "E7" "1.65%"
Set-TextValue "G7" "13"
Set-TextValue "D8" "3.725"
Set-TextValue "E8" "0.11%"
Set-TextValue "G8" "13"
Set-TextValue "D9" "0.9129"
Set-TextValue "E9" "1.50%"
Set-TextValue "G9" "13"
Set-TextValue "E10" "17.42%"
Set-TextValue "G10" "13"
Set-TextValue "D11" "0.1693"
Set-TextValue "E11" "1.33%"
Set-TextValue "G11" "13"
Set-TextValue "D12" "0.08231"
Set-TextValue "E12" "2.65%"
Set-TextValue "G12" "13"
Set-TextValue "D13" "0.03119"
Set-TextValue "E13" "2.59%"
Set-TextValue "G13" "13"
Set-TextValue "D14" "0.09959"
Set-TextValue "E14" "-0.47%"
Set-TextValue "G14" "13"
Set-TextValue "D15" "0.001494"
Set-TextValue "E15" "0.03%"
Set-TextValue "G15" "13"
Set-TextValue "D16" "0.005756"
Set-TextValue "E16" "1.05%"
Set-TextValue "G16" "13"
Set-TextValue "D17" "3.499"
Set-TextValue "E17" "0.37%"
Set-TextValue "G17" "13"
Set-TextValue "D18" "2.098"
Set-TextValue "E18" "1.08%"
Set-TextValue "G18" "13"
Set-TextValue "D19" "0.3324"
Set-TextValue "E19" "0.70%"
Set-TextValue "G19" "13"
Set-TextValue "D20" "0.1287"
Set-TextValue "E20" "-0.62%"
Set-TextValue "G20" "13"
Set-TextValue "D21" "4.155"
Set-TextValue "E21" "4.13%"
Set-TextValue "G21" "13"
Set-TextValue "D22" "0.2098"
Set-TextValue "E22" "-0.25%"
Set-TextValue "G22" "13"
Set-TextValue "D23" "0.04518"
Set-TextValue "E23" "0.44%"
Set-TextValue "G23" "13"
Set-TextValue "D24" "0.001206"
Set-TextValue "E24" "-0.77%"
Set-TextValue "G24" "13"
Set-TextValue "D25" "0.004183"
Set-TextValue "E25" "-9.42%"
Set-TextValue "G25" "13"
Set-TextValue "D26" "0.0001298"
Set-TextValue "E26" "-0.23%"
Set-TextValue "G26" "13"
Set-TextValue "D27" "0.0003388"
Set-TextValue "G27" "13"
Set-TextValue "G28" "13"
Set-TextValue "G29" "13"
Set-TextValue "G30" "13"
Set-TextValue "G31" "13"
Set-TextValue "G32" "13"
Set-TextValue "G33" "13"
Set-TextValue "G34" "13"
Set-TextValue "G35" "13"
Set-TextValue "G36" "13"
Set-TextValue "G37" "13"
Set-TextValue "G38" "13"
Set-TextValue "D39" "0.01579"
Set-TextValue "E39" "1.00%"
Set-TextValue "G39" "13"
Set-TextValue "D40" "0.04477"
Set-TextValue "E40" "3.53%"
Set-TextValue "G40" "13"
Set-TextValue "D41" "0.007379"
Set-TextValue "E41" "0.28%"
Set-TextValue "G41" "13"
Set-TextValue "D42" "0.009505"
Set-TextValue "E42" "-5.48%"
Set-TextValue "G42" "13"
Set-TextValue "D43" "0.1333"
Set-TextValue "E43" "2.58%"
Set-TextValue "G43" "13"
Set-TextValue "D44" "0.002237"
Set-TextValue "E44" "10.28%"
Set-TextValue "G44" "13"
Set-TextValue "D45" "0.008018"
Set-TextValue "E45" "-12.06%"
Set-TextValue "G45" "13"
Set-TextValue "D46" "0.00006094"
Set-TextValue "E46" "3.12%"
Set-TextValue "G46" "13"
Set-TextValue "E47" "-0.29%"
Set-TextValue "G47" "13"
Set-TextValue "D48" "2.606"
Set-TextValue "E48" "15.38%"
Set-TextValue "G48" "13"
Set-TextValue "D49" "0.001996"
Set-TextValue "E49" "-33.45%"
Set-TextValue "G49" "13"
Set-TextValue "D50" "0.00002096"
Set-TextValue "E50" "-0.29%"
Set-TextValue "G50" "13"
Set-TextValue "D51" "0.0001996"
Set-TextValue "E51" "-0.29%"
Set-TextValue "G51" "13"
